$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet originally held 6 rows of a single comma/quote-joined "CSV" string
# in column A (e.g. `Name,"Age","Occupation","Country"`). This splits that
# data out into four real columns (A:D) the way Excel's "Convert Text to
# Columns" wizard would, with column A kept/re-applied as Text (@) format and
# the numeric Age column stored as real numbers.
# ---------------------------------------------------------------------------

# Row 1 - header
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Age"
$ws.Range("C1").Value = "Occupation"
$ws.Range("D1").Value = "Country"

# Row 2
$ws.Range("A2").Value = "John Doe"
$ws.Range("B2").Value = 28
$ws.Range("C2").Value = "Software Developer"
$ws.Range("D2").Value = "USA"

# Row 3
$ws.Range("A3").Value = "Jane Smith"
$ws.Range("B3").Value = 34
$ws.Range("C3").Value = "Architect"
$ws.Range("D3").Value = "Canada"

# Row 4
$ws.Range("A4").Value = "Robert Brown"
$ws.Range("B4").Value = 22
$ws.Range("C4").Value = "Student"
$ws.Range("D4").Value = "UK"

# Row 5
$ws.Range("A5").Value = "Maria Garcia"
$ws.Range("B5").Value = 45
$ws.Range("C5").Value = "Teacher"
$ws.Range("D5").Value = "Spain"

# Row 6
$ws.Range("A6").Value = "Xiu Ying"
$ws.Range("B6").Value = 30
$ws.Range("C6").Value = "Engineer"
$ws.Range("D6").Value = "China"

# Column A keeps a Text number format (matches cellXfs xf numFmtId="49")
$ws.Range("A1:A6").NumberFormat = "@"

# Column A is left sized to fit the original (wider, pre-split) text.
$ws.Columns("A").ColumnWidth = 36.71

# Leftover UI state from the wizard: active cell on A5.
$ws.Range("A5").Select() | Out-Null
